$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'92.545.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.91%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.110.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.40%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'235.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.74%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'613.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.99%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.09"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.02%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.389"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.40%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -0.02%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'3.109.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.41%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.784"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.19%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.198"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.28%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -3.72%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'92.177.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.90%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'33.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.08%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  -3.49%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.695.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.31%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.127.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.02%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'3.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.63%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'14.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.59%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'5.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.24%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.0000203"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.29%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'9.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.47%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'439.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.27%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  -6.17%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'85.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.45%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'11.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.78%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'3.274.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.53%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.11%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.179"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +6.47%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.237"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.50%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.123"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -16.76%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'9.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.64%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -35.23%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'8.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +8.18%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.159"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -8.83%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'25.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.51%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'3.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.41%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'1.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.88%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +7.69%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -2.75%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'467.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.00%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -1.84%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'3.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.10%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  +0.06%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'160.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.46%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.683"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.30%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -4.76%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "'1.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.05%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0328"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.04%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'43.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.43%  "
$ws.Range("E51").Style = "Normal"
